$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three blank spacer rows that used to sit between the title block
# (rows 1:5) and the data table (old rows 9:12) were removed, so the
# data table now starts right after the title.
$ws.Range("A6:A8").EntireRow.Delete()

# After the shift, the data table occupies rows 6:9.
# Row 6 = headers (Cong viet / Ngay / Nhiem vu / Decline) - unchanged.
# Row 7 used to be split across two rows/cells ("Gap giao vien va nghe
# huong dan " in col A + "ve de tai da nhan" on the row below); it is
# now one corrected, combined sentence in a single cell.
$ws.Range("A7").Value = "Gặp giáo viên và nghe hướng dẫn về đề tài đã nhận"
$ws.Range("A8").ClearContents()

# The "2. Ve giao dien..." task description was reworded/shortened and
# the typo "dien" -> "dien" (missing dau nang) was fixed.
$ws.Range("C8").Value = "2. Vẽ giao diện cho app (cho User và Doctor)"

# Column A was widened (it now holds the longer combined sentence) and
# column E was widened as well.
$ws.Columns("A").ColumnWidth = 45.28515625
$ws.Columns("E").ColumnWidth = 14.42578125

# Restore the cursor/selection to where it was left in the saved file.
[void]$ws.Range("C17").Select()
